$wb = $excel.ActiveWorkbook

# --- Restrictions sheet (sheet6): add a NumIterations column for the
#     upcoming Monte Carlo simulation work ---
$wsRestrictions = $wb.Worksheets.Item("Restrictions")
$wsRestrictions.Range("B1").Value = "NumIterations"
$wsRestrictions.Range("B2").Value = 10

# Match the original workbook's "bestFit" column sizing as closely as this
# runtime's column-width model allows.
$wsRestrictions.Columns.Item(1).ColumnWidth = 7.75
$wsRestrictions.Columns.Item(2).ColumnWidth = 13.17

# --- Selection / active-sheet bookkeeping ---
# Sources: just move the cursor, sheet stays non-active.
$wsSources = $wb.Worksheets.Item("Sources")
$wsSources.Range("I4").Select()

# Restrictions becomes the active / selected tab (previously it was Hubs).
$wsRestrictions.Range("B2").Select()
